# Adds 45X for grid batteries; adjusts 45X implementation in transport sector
#
# 1. Insert a new "Phase Out" worksheet between "About" and "BSfVBP"
#    (sheetId 3, placed right after About) holding the yearly phase-out
#    factors for 2023-2032.
# 2. Re-point the BSfVBP sheet's yearly price formulas (F2:M2) so they pull
#    their phase-out multiplier from the new "Phase Out" sheet, and compound
#    the annual inflation-adjustment directly off the base $E2 price (rather
#    than chaining off the previous year's cell) using the year offset from
#    $E1.

$wb = $excel.ActiveWorkbook

# --- 1. Add "Phase Out" worksheet right after "About" -----------------
$aboutSheet = $wb.Worksheets.Item("About")
$phaseOut = $wb.Worksheets.Add($null, $aboutSheet)
$phaseOut.Name = "Phase Out"

$years   = @(2023, 2024, 2025, 2026, 2027, 2028, 2029, 2030, 2031, 2032)
$factors = @(1,    1,    1,    1,    1,    1,    1,    0.75, 0.5,  0.25)
$cols    = @("B",  "C",  "D",  "E",  "F",  "G",  "H",  "I",  "J",  "K")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $phaseOut.Range($cols[$i] + "2").Value = $years[$i]
    $phaseOut.Range($cols[$i] + "3").Value = $factors[$i]
}

$phaseOut.Range("L3").Select() | Out-Null

# --- 2. Update the BSfVBP price formulas (F2:M2) -----------------------
$bsf = $wb.Worksheets.Item("BSfVBP")

$targetCols = @("F", "G", "H", "I", "J", "K", "L", "M")
$phaseCols  = @("D", "E", "F", "G", "H", "I", "J", "K")

for ($i = 0; $i -lt $targetCols.Length; $i++) {
    $tc = $targetCols[$i]
    $pc = $phaseCols[$i]
    $bsf.Range($tc + "2").Formula = "=`$E2*(1-About!`$A`$14)^(" + $tc + "1-`$E1)*'Phase Out'!" + $pc + "3"
}

$bsf.Range("D36").Select() | Out-Null

# --- 3. Restore "About" sheet's prior selection -------------------------
$aboutSheet.Range("B37").Select() | Out-Null
